$d = $word.ActiveDocument

$replacements = @(
    @("677×2=", "684×5="),
    @("178×2=", "773×3="),
    @("212×7=", "665×7="),
    @("460×9=", "435×4="),
    @("699×6=", "563×9="),
    @("458×6=", "930×6="),
    @("707×5=", "815×9="),
    @("850×4=", "461×2="),
    @("690×5=", "255×6="),
    @("930×9=", "893×5="),
    @("635×2=", "658×9="),
    @("793×3=", "177×2="),
    @("535×2=", "457×4="),
    @("337×6=", "635×3="),
    @("764×7=", "950×7="),
    @("655×7=", "720×3="),
    @("426×2=", "888×8="),
    @("537×8=", "381×8="),
    @("319×4=", "674×5="),
    @("546×6=", "386×7="),
    @("512×9=", "741×6="),
    @("836×4=", "824×7="),
    @("669×5=", "420×3="),
    @("813×8=", "450×9="),
    @("313×4=", "469×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
